# edit.ps1 - applies the "Add files via upload" diff to before.docx
$d = $word.ActiveDocument

function Merge-ParaRange($doc, $startPos, $text) {
    $len = $text.Length
    $r = $doc.Range($startPos, $startPos + $len)
    $r.Text = "#"
    $r2 = $doc.Range($startPos, $startPos + 1)
    $r2.Text = $text
}

$p = $d.Paragraphs.Item(4)
Merge-ParaRange $d $p.Range.Start "Dans page admin, 4 lien, ajout participants, création, surpression et modif event."

$p = $d.Paragraphs.Item(6)
Merge-ParaRange $d $p.Range.Start "Ajout de participant : nom – prenom – email – date_naiss – organisation/entreprise – observations – event auquelle il participe ( un seul event ici ) "

$p = $d.Paragraphs.Item(7)
Merge-ParaRange $d $p.Range.Start "Création de event : titre – thème – date_debut – durée en jours – nb max participant – description – organisateur (nom) – type d-event "

$p = $d.Paragraphs.Item(11)
Merge-ParaRange $d $p.Range.Start "Si ajout de même adresse mail à un event : on fusionne "

$p = $d.Paragraphs.Item(12)
Merge-ParaRange $d $p.Range.Start "Si nb max de participant à un event : bloque "

$p = $d.Paragraphs.Item(17)
Merge-ParaRange $d $p.Range.Start "table event, lien entre event et participant "

# Paragraph 5: protect run 6 (" et ajout participant") from forward merge-cascade
$p5 = $d.Paragraphs.Item(5)
$p5Start = $p5.Range.Start
$p5rng = $p5.Range
$found5 = $p5rng.Find.Execute(" et ajout participant")
$p5rng.Font.Size = 11
Merge-ParaRange $d $p5Start "Pages creation, modif et suppression"

# Paragraph 15: protect the single-space run (offset 11) from backward merge-cascade
$p15 = $d.Paragraphs.Item(15)
$p15Start = $p15.Range.Start
$p15SpaceRng = $d.Range($p15Start + 11, $p15Start + 12)
$p15SpaceRng.Font.Size = 11
Merge-ParaRange $d ($p15Start + 12) "avec nom et mdp correspondant "

# Append an empty paragraph and a new paragraph with the project URL at the end
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertParagraphAfter()
$endRange2 = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange2.InsertParagraphAfter()
$urlRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$urlRange.InsertAfter("http://cboin.hd.free.fr/jee/")

Write-Output "done"
